# The workbook's text cells contain a mangled "plus-minus" sign: the UTF-8
# bytes for U+00B1 (±) were themselves decoded as Latin-1/Windows-1252 and
# re-encoded, producing the two-character mojibake sequence U+00C2 U+00B1
# ("Â±"). This script restores the correct single "±" character everywhere
# it occurs, without otherwise touching any cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badSeq  = [string][char]0x00C2 + [string][char]0x00B1   # "Â±"
$goodSeq = [string][char]0x00B1                            # "±"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$rowStart = $used.Row
$colStart = $used.Column

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($rowStart + $r, $colStart + $c)
        $val = $cell.Value2()
        if ($val -ne $null -and $val -is [string] -and $val.Contains($badSeq)) {
            $cell.Value2 = $val.Replace($badSeq, $goodSeq)
        }
    }
}
